# Remove the "「クウェート国」" entry (row 69) from the posts sheet.
# This shifts every subsequent row up by one and shrinks the used range
# from A1:C242 to A1:C241.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(69).Delete()
